$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("C2").Value = 87
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 3379.408121109009
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 15

# Update row 3
$ws.Range("C3").Value = 89
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 2532.006978988647
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 15

# Update row 4
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 2782.976865768433
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 15

# Delete row 5 entirely (shift cells up)
$ws.Range("A5:G5").Delete()
